$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; this shifts the existing rows 10-62
# down to 11-63 (preserving all of their data/styles), matching the
# diff's "new_row[N] = old_row[N-1]" shift pattern plus the dimension
# growing from A1:T62 to A1:T63.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with a new weekly entry. It mirrors
# the (now shifted-down) original row 10 except for the date (D) and
# volume (M), which carry the new reported values.
$ws.Cells.Item(10, 1).Value = 5
$ws.Cells.Item(10, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(10, 3).Value = "Maule"
$ws.Cells.Item(10, 4).Value = 44622
$ws.Cells.Item(10, 5).Value = 7
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100101
$ws.Cells.Item(10, 8).Value = "Berries"
$ws.Cells.Item(10, 9).Value = 100101001
$ws.Cells.Item(10, 10).Value = "Arándano (blue)"
$ws.Cells.Item(10, 11).Value = "Sin especificar"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 80
$ws.Cells.Item(10, 14).Value = 3000
$ws.Cells.Item(10, 15).Value = 3000
$ws.Cells.Item(10, 16).Value = 3000
$ws.Cells.Item(10, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(10, 18).Value = "Provincia de Linares"
$ws.Cells.Item(10, 19).Value = 1500
$ws.Cells.Item(10, 20).Value = 2
